$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7943037748336792
$ws.Range("B1").Value = 2.093758821487427
$ws.Range("D1").Value = 1.250910758972168
$ws.Range("E1").Value = 0.5209749937057495
